$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '27.525.56'
$ws.Range("E2").Value = '  +1.76%  '

# Row 3
$ws.Range("D3").Value = '1.567.18'
$ws.Range("E3").Value = '  +0.14%  '

# Row 4
$ws.Range("E4").Value = '  -1.53%  '

# Row 5
$ws.Range("D5").Value = '''210.95'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +1.18%  '

# Row 6
$ws.Range("E6").Value = '  -0.37%  '

# Row 7
$ws.Range("D7").Value = '''0.989'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -1.57%  '

# Row 8
$ws.Range("D8").Value = '''22.61'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +2.45%  '

# Row 9
$ws.Range("E9").Value = '  +0.37%  '

# Row 10
$ws.Range("E10").Value = '  -0.45%  '

# Row 11
$ws.Range("D11").Value = '''0.0869'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +1.16%  '

# Row 12
$ws.Range("D12").Value = '1.791.39'
$ws.Range("E12").Value = '  +0.20%  '

# Row 13
$ws.Range("D13").Value = '1.567.52'
$ws.Range("E13").Value = '  +0.09%  '

# Row 14
$ws.Range("E14").Value = '  +0.07%  '

# Row 15
$ws.Range("E15").Value = '  -0.27%  '

# Row 16
$ws.Range("D16").Value = '27.512.68'
$ws.Range("E16").Value = '  +1.79%  '

# Row 17
$ws.Range("D17").Value = '''62.40'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +0.80%  '

# Row 18
$ws.Range("D18").Value = '''225.87'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +4.49%  '

# Row 19
$ws.Range("E19").Value = '  +1.33%  '

# Row 20
$ws.Range("D20").Value = '0.0₃0704'
$ws.Range("E20").Value = '  -0.12%  '

# Row 21
$ws.Range("E21").Value = '  -1.56%  '

# Row 22
$ws.Range("E22").Value = '  -0.61%  '

# Row 23
$ws.Range("D23").Value = '''9.43'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +2.34%  '

# Row 24
$ws.Range("D24").Value = '''1.95'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.04%  '

# Row 25
$ws.Range("D25").Value = '''149.95'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -2.58%  '

# Row 26
$ws.Range("B26").Value = 'EthereumClassic'
$ws.Range("C26").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D26").Value = '''15.17'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +0.60%  '

# Row 27
$ws.Range("B27").Value = 'Cosmos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D27").Value = '''6.61'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -0.10%  '

# Row 28
$ws.Range("E28").Value = '  +1.62%  '

# Row 29
$ws.Range("E29").Value = '  -1.53%  '

# Row 30
$ws.Range("E30").Value = '  +0.91%  '

# Row 31
$ws.Range("E31").Value = '  -0.60%  '

# Row 32
$ws.Range("D32").Value = '''3.24'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +0.13%  '

# Row 33
$ws.Range("D33").Value = '1.447.06'
$ws.Range("E33").Value = '  +1.23%  '

# Row 34
$ws.Range("E34").Value = '  -1.92%  '

# Row 35
$ws.Range("D35").Value = '''1.11'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +2.35%  '

# Row 36
$ws.Range("D36").Value = '''1.60'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -0.22%  '

# Row 37
$ws.Range("E37").Value = '  -0.81%  '

# Row 38
$ws.Range("D38").Value = '''0.0168'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +0.26%  '

# Row 39
$ws.Range("D39").Value = '''0.541'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +1.43%  '

# Row 40
$ws.Range("E40").Value = '  -0.03%  '

# Row 41
$ws.Range("B41").Value = 'MXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D41").Value = '''2.37'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +1.82%  '

# Row 42
$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").Value = '''5.72'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -1.18%  '

# Row 43
$ws.Range("D43").Value = '''0.989'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -1.56%  '

# Row 44
$ws.Range("D44").Value = '''1.84'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +5.42%  '

# Row 45
$ws.Range("E45").Value = '  -3.00%  '

# Row 46
$ws.Range("D46").Value = '''64.59'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -0.30%  '

# Row 47
$ws.Range("D47").Value = '1.703.50'
$ws.Range("E47").Value = '  +0.31%  '

# Row 48
$ws.Range("E48").Value = '  -0.01%  '

# Row 49
$ws.Range("E49").Value = '  +1.24%  '

# Row 50
$ws.Range("D50").Value = '''0.0946'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -1.52%  '

# Row 51
$ws.Range("E51").Value = '  -1.52%  '
